$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row labelled "Docentes responsáveis:" (row 12) is removed entirely; every
# row below it shifts up by one (this also carries the correct row heights and
# most B/C cell content up with it).
$ws.Rows.Item(12).Delete()

# A few rows need their B/C (value + highlighted-value) text corrected after the
# shift, because the original sheet had mis-aligned label/value pairs.

# Row 10 "Objetivos:" now gets the real (Portuguese) objectives text.
$ws.Range("B10").Value = "A disciplina visa apresentar aos estudantes o conceito, tipos, modelos e sistemas de eco-inovação para o desenvolvimento da capacidade analítica e propositiva como competências profissionais nas áreas de inovação e sustentabilidade."
$ws.Range("C10").Value = "A disciplina visa apresentar aos estudantes o conceito, tipos, modelos e sistemas de eco-inovação para o desenvolvimento da capacidade analítica e propositiva como competências profissionais nas áreas de inovação e sustentabilidade."

# Row 12 "Programa resumido:" now carries the professor identification text.
$ws.Range("B12").Value = "5840820 - Gustavo Aristides Santana Martinez"
$ws.Range("C12").Value = "5840820 - Gustavo Aristides Santana Martinez"

# Row 14 "Programa:" now gets the (Portuguese) short program description.
$ws.Range("B14").Value = "Eco inovação. Métricas da eco-inovação. Introdução ao Ciclo de vida do produto. Eco inovação na indústria. Estudo de casos de projetos de eco-inovação no Brasil. Métodos e ferramentas suporte do processo de eco-inovação. Identificação antecipada de falha como suporte a eco-inovação. TRIZ como resposta a eco-inovação. Proposta metodológica para soluções eco-inovadoras."
$ws.Range("C14").Value = "Eco inovação. Métricas da eco-inovação. Introdução ao Ciclo de vida do produto. Eco inovação na indústria. Estudo de casos de projetos de eco-inovação no Brasil. Métodos e ferramentas suporte do processo de eco-inovação. Identificação antecipada de falha como suporte a eco-inovação. TRIZ como resposta a eco-inovação. Proposta metodológica para soluções eco-inovadoras."

# Row 17 "Método:" now gets the full (Portuguese) detailed program/method text.
$ws.Range("B17").Value = "1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações.2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD.3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management.4. Eco inovação na indústria: química, agro alimentos, metal mecânica.5. Estudo de casos de projetos de eco inovação no Brasil.6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros.7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos.8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições.9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar"
$ws.Range("C17").Value = "1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações.2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD.3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management.4. Eco inovação na indústria: química, agro alimentos, metal mecânica.5. Estudo de casos de projetos de eco inovação no Brasil.6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros.7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos.8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições.9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar"
